$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data corrections (Sunidhi -> Simran / Mahalinga) ---
$ws.Range("C2").Value = "Simran"

# --- New columns V, W, X: headers (row 1) and data (row 2) ---
# Column V (22): Branch / Indore I
$ws.Range("V1").Value = "Branch"
$ws.Range("S1").Copy()
$ws.Range("V1").PasteSpecial(-4122)

$ws.Range("V2").Value = "Indore I"
$ws.Range("S2").Copy()
$ws.Range("V2").PasteSpecial(-4122)

# Column W (23): UserManagementPageRole / Admin
$ws.Range("W1").Value = "UserManagementPageRole"
$ws.Range("S1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

$ws.Range("W2").Value = "Admin"
$ws.Range("S2").Copy()
$ws.Range("W2").PasteSpecial(-4122)

# Column X (24): UsernameInUNPage / IBU0001192
$ws.Range("X1").Value = "UsernameInUNPage"
$ws.Range("S1").Copy()
$ws.Range("X1").PasteSpecial(-4122)

$ws.Range("X2").Value = "IBU0001192"
$ws.Range("S2").Copy()
$ws.Range("X2").PasteSpecial(-4122)

# --- Remaining row2 data corrections ---
$ws.Range("P2").Value = "Mahalinga"
$ws.Range("D2").Value = "mahalinga@gmail.com"

# --- New column widths ---
$ws.Columns.Item(22).ColumnWidth = 11.02
$ws.Columns.Item(23).ColumnWidth = 25.45
$ws.Columns.Item(24).ColumnWidth = 17.59

# --- Sheet view: clear frozen/topLeftCell scroll position, update selection ---
$ws.Range("D7").Select()

Write-Host "edit applied"
